$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "prov_full" column (D) giving the full official name for each
# province abbreviation already listed in column A. Cells are written in the
# same order the author's workbook accumulated them in (matches the shared
# string table insertion order), then the header is written last.
$ws.Cells.Item(6, 4).Value = "安徽省"
$ws.Cells.Item(13, 4).Value = "云南省"
$ws.Cells.Item(4, 4).Value = "湖南省"
$ws.Cells.Item(17, 4).Value = "广西壮族自治区"
$ws.Cells.Item(18, 4).Value = "河北省"
$ws.Cells.Item(26, 4).Value = "陕西省"
$ws.Cells.Item(31, 4).Value = "青海省"
$ws.Cells.Item(5, 4).Value = "湖北省"
$ws.Cells.Item(24, 4).Value = "甘肃省"
$ws.Cells.Item(21, 4).Value = "新疆维吾尔自治区"
$ws.Cells.Item(25, 4).Value = "辽宁省"
$ws.Cells.Item(29, 4).Value = "吉林省"
$ws.Cells.Item(22, 4).Value = "河南省"
$ws.Cells.Item(11, 4).Value = "四川省"
$ws.Cells.Item(12, 4).Value = "江西省"
$ws.Cells.Item(14, 4).Value = "贵州省"
$ws.Cells.Item(27, 4).Value = "山西省"
$ws.Cells.Item(28, 4).Value = "内蒙古自治区"
$ws.Cells.Item(23, 4).Value = "黑龙江省"
$ws.Cells.Item(3, 4).Value = "广东省"
$ws.Cells.Item(16, 4).Value = "福建省"
$ws.Cells.Item(32, 4).Value = "西藏自治区"
$ws.Cells.Item(20, 4).Value = "海南省"
$ws.Cells.Item(2, 4).Value = "浙江省"
$ws.Cells.Item(30, 4).Value = "宁夏回族自治区"
$ws.Cells.Item(19, 4).Value = "江苏省"
$ws.Cells.Item(15, 4).Value = "山东省"
$ws.Cells.Item(7, 4).Value = "天津市"
$ws.Cells.Item(8, 4).Value = "北京市"
$ws.Cells.Item(9, 4).Value = "上海市"
$ws.Cells.Item(10, 4).Value = "重庆市"
$ws.Cells.Item(1, 4).Value = "prov_full"

# F12 picked up a stray leading apostrophe (text quote-prefix) while editing
# and was then cleared, leaving an empty "quote prefix" formatted cell.
$ws.Range("F12").Value = "'"
$ws.Range("F12").Value = ""

# Window/view tweaks made while reviewing the new column.
$excel.ActiveWindow.Zoom = 108
$null = $ws.Range("D30").Select()
